# Updated cryptos list on Mon Jul  3 23:43:33 UTC 2023 with GitHub Actions
# Refresh prices / 1h volume % for all listed coins; a new coin (OKB) was
# inserted at rank 7, shifting every subsequent coin down one row and
# dropping the former last row (Elrond) off the bottom of the table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '31.110.64'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.76%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.953.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.01%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9997'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '245.58'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.02%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9994'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.07%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4889'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.97%  '

# Row 8
$ws.Range("B8").Value = 'OKB'
$ws.Range("C8").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '44.76'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.47%  '

# Row 9
$ws.Range("B9").Value = 'Cardano'
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2958'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.56%  '

# Row 10
$ws.Range("B10").Value = 'Dogecoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06821'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.43%  '

# Row 11
$ws.Range("B11").Value = 'Solana'
$ws.Range("C11").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.29%  '

# Row 12
$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '106.57'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -6.18%  '

# Row 13
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07725'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.73%  '

# Row 14
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.927.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.41%  '

# Row 15
$ws.Range("B15").Value = 'Polkadot'
$ws.Range("C15").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.407'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.49%  '

# Row 16
$ws.Range("B16").Value = 'Polygon'
$ws.Range("C16").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.7119'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +4.75%  '

# Row 17
$ws.Range("B17").Value = 'BitcoinCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '285.77'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -4.28%  '

# Row 18
$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '30.993.25'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.34%  '

# Row 19
$ws.Range("B19").Value = 'ShibaInu'
$ws.Range("C19").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007741'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.27%  '

# Row 20
$ws.Range("B20").Value = 'Avalanche'
$ws.Range("C20").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.19'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.75%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9997'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.05%  '

# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.180.04'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.28%  '

# Row 23
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.506'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '

# Row 24
$ws.Range("B24").Value = 'BinanceUSD'
$ws.Range("C24").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9999'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.03%  '

# Row 25
$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '6.600'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.33%  '

# Row 26
$ws.Range("B26").Value = 'Cosmos'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.909'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.58%  '

# Row 27
$ws.Range("B27").Value = 'Monero'
$ws.Range("C27").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '168.80'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.23%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.97'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.03%  '

# Row 29
$ws.Range("B29").Value = 'LidoDAOToken'
$ws.Range("C29").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.192'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.14%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.1051'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -1.71%  '

# Row 31
$ws.Range("B31").Value = 'Toncoin'
$ws.Range("C31").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.438'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +0.52%  '

# Row 32
$ws.Range("B32").Value = 'Filecoin'
$ws.Range("C32").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.741'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +15.73%  '

# Row 33
$ws.Range("B33").Value = 'InternetComputer(DFINITY)'
$ws.Range("C33").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.467'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.93%  '

# Row 34
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05000'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.03%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7604'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.33%  '

# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.164'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +1.44%  '

# Row 37
$ws.Range("B37").Value = 'HuobiToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.728'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.41%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02040'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.14%  '

# Row 39
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.704'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.41%  '

# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.148'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +5.79%  '

# Row 41
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.430'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +10.30%  '

# Row 42
$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4484'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.48%  '

# Row 43
$ws.Range("B43").Value = 'Quant'
$ws.Range("C43").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '109.46'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.68%  '

# Row 44
$ws.Range("B44").Value = 'TrustWalletToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.8789'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.78%  '

# Row 45
$ws.Range("B45").Value = 'Aave'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '72.41'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +3.75%  '

# Row 46
$ws.Range("B46").Value = 'PaxDollar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.9986'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.31%  '

# Row 47
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.474'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.14%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '979.08'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +15.57%  '

# Row 49
$ws.Range("B49").Value = 'Algorand'
$ws.Range("C49").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.1281'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.85%  '

# Row 50
$ws.Range("B50").Value = 'EnergySwap'
$ws.Range("C50").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '9.323'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.23%  '

# Row 51
$ws.Range("B51").Value = 'WOONetwork'
$ws.Range("C51").Value = 'https://coinranking.com/coin/k-J3YwacF+woonetwork-woo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.2594'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.83%  '
